$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column F (Tamany Lliure GB) values ---
$ws.Range("F3:F7").Value = 1764
$ws.Range("F8:F9").Value = 1105
$ws.Range("F10:F29").Value = 1646
$ws.Range("F30:F35").Value = 4748
$ws.Range("F36:F37").Value = 3157
$ws.Range("F38:F39").Value = 426
$ws.Range("F40:F66").Value = 7177

# --- Update column A (Nom NAS) values ---
$ws.Range("A38:A39").Value = "TorinoLocal"
$ws.Range("A40:A66").Value = "TransportsTresserras"

# --- Update column B (Nom Dispositiu) values ---
$ws.Range("B36").Value = "SIS0010SV0721"
$ws.Range("B37").Value = "SIS0011SV0721"
$ws.Range("B38").Value = "SERVER1"
$ws.Range("B39").Value = "PCDIRECCIO"
$ws.Range("B40").Value = "Server235 - TS B1 Nou"
$ws.Range("B41").Value = "Server226"
$ws.Range("B42").Value = "PcVirtual1 - planella"
$ws.Range("B43").Value = "Server215 - KSC"
$ws.Range("B44").Value = "Server219 - TS"
$ws.Range("B45").Value = "Server220 - Web"
$ws.Range("B46").Value = "Server236 - Qlik"
$ws.Range("B47").Value = "pcvirutal09_win10 - EDI"
$ws.Range("B48").Value = "SERVER203 - Progress"
$ws.Range("B49").Value = "Server221- Aduana"
$ws.Range("B50").Value = "Server223"
$ws.Range("B51").Value = "Server224 - AD2k12 - restored"
$ws.Range("B52").Value = "Server228 - GMagat"
$ws.Range("B53").Value = "Server229 - IP6"
$ws.Range("B54").Value = "Server238 - InterCompany"
$ws.Range("B55").Value = "Server244 - Dades"
$ws.Range("B56").Value = "pcvirtual08_win7 - Spiceworks"
$ws.Range("B57").Value = "Server227"
$ws.Range("B58").Value = "ServerBCN220 - dcserver BCN"
$ws.Range("B59").Value = "SERVER208"
$ws.Range("B60").Value = "SERVER222"
$ws.Range("B61").Value = "Server200-VCenter"
$ws.Range("B62").Value = "server251 - Unifi controller"
$ws.Range("B63").Value = "Server237 - SAPupdate (1)"
$ws.Range("B64").Value = "Server226"
$ws.Range("B65").Value = "Server227"
$ws.Range("B66").Value = "Server235 - TS B1 Nou"

# --- Update formulas in columns C/D/E for restructured rows ---
$ws.Range("C38").Formula = "=LOOKUP(2,1/(TorinoLocal!1:1<>""""),TorinoLocal!1:1)"
$ws.Range("D38").Formula = "=LOOKUP(2,1/(TorinoLocal!3:3<>""""),TorinoLocal!3:3)"
$ws.Range("E38").Formula = "=LOOKUP(2,1/(TorinoLocal!2:2<>""""),TorinoLocal!2:2)"
$ws.Range("C39").Formula = "=LOOKUP(2,1/(TorinoLocal!5:5<>""""),TorinoLocal!5:5)"
$ws.Range("D39").Formula = "=LOOKUP(2,1/(TorinoLocal!7:7<>""""),TorinoLocal!7:7)"
$ws.Range("E39").Formula = "=LOOKUP(2,1/(TorinoLocal!6:6<>""""),TorinoLocal!6:6)"
$ws.Range("C40").Formula = "=LOOKUP(2,1/(TransportsTresserras!1:1<>""""),TransportsTresserras!1:1)"
$ws.Range("D40").Formula = "=LOOKUP(2,1/(TransportsTresserras!3:3<>""""),TransportsTresserras!3:3)"
$ws.Range("E40").Formula = "=LOOKUP(2,1/(TransportsTresserras!2:2<>""""),TransportsTresserras!2:2)"
$ws.Range("C41").Formula = "=LOOKUP(2,1/(TransportsTresserras!5:5<>""""),TransportsTresserras!5:5)"
$ws.Range("D41").Formula = "=LOOKUP(2,1/(TransportsTresserras!7:7<>""""),TransportsTresserras!7:7)"
$ws.Range("E41").Formula = "=LOOKUP(2,1/(TransportsTresserras!6:6<>""""),TransportsTresserras!6:6)"
$ws.Range("C42").Formula = "=LOOKUP(2,1/(TransportsTresserras!9:9<>""""),TransportsTresserras!9:9)"
$ws.Range("D42").Formula = "=LOOKUP(2,1/(TransportsTresserras!11:11<>""""),TransportsTresserras!11:11)"
$ws.Range("E42").Formula = "=LOOKUP(2,1/(TransportsTresserras!10:10<>""""),TransportsTresserras!10:10)"
$ws.Range("C43").Formula = "=LOOKUP(2,1/(TransportsTresserras!13:13<>""""),TransportsTresserras!13:13)"
$ws.Range("D43").Formula = "=LOOKUP(2,1/(TransportsTresserras!15:15<>""""),TransportsTresserras!15:15)"
$ws.Range("E43").Formula = "=LOOKUP(2,1/(TransportsTresserras!14:14<>""""),TransportsTresserras!14:14)"
$ws.Range("C44").Formula = "=LOOKUP(2,1/(TransportsTresserras!17:17<>""""),TransportsTresserras!17:17)"
$ws.Range("D44").Formula = "=LOOKUP(2,1/(TransportsTresserras!19:19<>""""),TransportsTresserras!19:19)"
$ws.Range("E44").Formula = "=LOOKUP(2,1/(TransportsTresserras!18:18<>""""),TransportsTresserras!18:18)"
$ws.Range("C45").Formula = "=LOOKUP(2,1/(TransportsTresserras!21:21<>""""),TransportsTresserras!21:21)"
$ws.Range("D45").Formula = "=LOOKUP(2,1/(TransportsTresserras!23:23<>""""),TransportsTresserras!23:23)"
$ws.Range("E45").Formula = "=LOOKUP(2,1/(TransportsTresserras!22:22<>""""),TransportsTresserras!22:22)"
$ws.Range("C46").Formula = "=LOOKUP(2,1/(TransportsTresserras!25:25<>""""),TransportsTresserras!25:25)"
$ws.Range("D46").Formula = "=LOOKUP(2,1/(TransportsTresserras!27:27<>""""),TransportsTresserras!27:27)"
$ws.Range("E46").Formula = "=LOOKUP(2,1/(TransportsTresserras!26:26<>""""),TransportsTresserras!26:26)"
$ws.Range("C47").Formula = "=LOOKUP(2,1/(TransportsTresserras!29:29<>""""),TransportsTresserras!29:29)"
$ws.Range("D47").Formula = "=LOOKUP(2,1/(TransportsTresserras!31:31<>""""),TransportsTresserras!31:31)"
$ws.Range("E47").Formula = "=LOOKUP(2,1/(TransportsTresserras!30:30<>""""),TransportsTresserras!30:30)"
$ws.Range("C48").Formula = "=LOOKUP(2,1/(TransportsTresserras!33:33<>""""),TransportsTresserras!33:33)"
$ws.Range("D48").Formula = "=LOOKUP(2,1/(TransportsTresserras!35:35<>""""),TransportsTresserras!35:35)"
$ws.Range("E48").Formula = "=LOOKUP(2,1/(TransportsTresserras!34:34<>""""),TransportsTresserras!34:34)"
$ws.Range("C49").Formula = "=LOOKUP(2,1/(TransportsTresserras!37:37<>""""),TransportsTresserras!37:37)"
$ws.Range("D49").Formula = "=LOOKUP(2,1/(TransportsTresserras!39:39<>""""),TransportsTresserras!39:39)"
$ws.Range("E49").Formula = "=LOOKUP(2,1/(TransportsTresserras!38:38<>""""),TransportsTresserras!38:38)"
$ws.Range("C50").Formula = "=LOOKUP(2,1/(TransportsTresserras!41:41<>""""),TransportsTresserras!41:41)"
$ws.Range("D50").Formula = "=LOOKUP(2,1/(TransportsTresserras!43:43<>""""),TransportsTresserras!43:43)"
$ws.Range("E50").Formula = "=LOOKUP(2,1/(TransportsTresserras!42:42<>""""),TransportsTresserras!42:42)"
$ws.Range("C51").Formula = "=LOOKUP(2,1/(TransportsTresserras!45:45<>""""),TransportsTresserras!45:45)"
$ws.Range("D51").Formula = "=LOOKUP(2,1/(TransportsTresserras!47:47<>""""),TransportsTresserras!47:47)"
$ws.Range("E51").Formula = "=LOOKUP(2,1/(TransportsTresserras!46:46<>""""),TransportsTresserras!46:46)"
$ws.Range("C52").Formula = "=LOOKUP(2,1/(TransportsTresserras!49:49<>""""),TransportsTresserras!49:49)"
$ws.Range("D52").Formula = "=LOOKUP(2,1/(TransportsTresserras!51:51<>""""),TransportsTresserras!51:51)"
$ws.Range("E52").Formula = "=LOOKUP(2,1/(TransportsTresserras!50:50<>""""),TransportsTresserras!50:50)"
$ws.Range("C53").Formula = "=LOOKUP(2,1/(TransportsTresserras!53:53<>""""),TransportsTresserras!53:53)"
$ws.Range("D53").Formula = "=LOOKUP(2,1/(TransportsTresserras!55:55<>""""),TransportsTresserras!55:55)"
$ws.Range("E53").Formula = "=LOOKUP(2,1/(TransportsTresserras!54:54<>""""),TransportsTresserras!54:54)"
$ws.Range("C54").Formula = "=LOOKUP(2,1/(TransportsTresserras!57:57<>""""),TransportsTresserras!57:57)"
$ws.Range("D54").Formula = "=LOOKUP(2,1/(TransportsTresserras!59:59<>""""),TransportsTresserras!59:59)"
$ws.Range("E54").Formula = "=LOOKUP(2,1/(TransportsTresserras!58:58<>""""),TransportsTresserras!58:58)"
$ws.Range("C55").Formula = "=LOOKUP(2,1/(TransportsTresserras!61:61<>""""),TransportsTresserras!61:61)"
$ws.Range("D55").Formula = "=LOOKUP(2,1/(TransportsTresserras!63:63<>""""),TransportsTresserras!63:63)"
$ws.Range("E55").Formula = "=LOOKUP(2,1/(TransportsTresserras!62:62<>""""),TransportsTresserras!62:62)"
$ws.Range("C56").Formula = "=LOOKUP(2,1/(TransportsTresserras!65:65<>""""),TransportsTresserras!65:65)"
$ws.Range("D56").Formula = "=LOOKUP(2,1/(TransportsTresserras!67:67<>""""),TransportsTresserras!67:67)"
$ws.Range("E56").Formula = "=LOOKUP(2,1/(TransportsTresserras!66:66<>""""),TransportsTresserras!66:66)"
$ws.Range("C57").Formula = "=LOOKUP(2,1/(TransportsTresserras!69:69<>""""),TransportsTresserras!69:69)"
$ws.Range("D57").Formula = "=LOOKUP(2,1/(TransportsTresserras!71:71<>""""),TransportsTresserras!71:71)"
$ws.Range("E57").Formula = "=LOOKUP(2,1/(TransportsTresserras!70:70<>""""),TransportsTresserras!70:70)"
$ws.Range("C58").Formula = "=LOOKUP(2,1/(TransportsTresserras!73:73<>""""),TransportsTresserras!73:73)"
$ws.Range("D58").Formula = "=LOOKUP(2,1/(TransportsTresserras!75:75<>""""),TransportsTresserras!75:75)"
$ws.Range("E58").Formula = "=LOOKUP(2,1/(TransportsTresserras!74:74<>""""),TransportsTresserras!74:74)"
$ws.Range("C59").Formula = "=LOOKUP(2,1/(TransportsTresserras!77:77<>""""),TransportsTresserras!77:77)"
$ws.Range("D59").Formula = "=LOOKUP(2,1/(TransportsTresserras!79:79<>""""),TransportsTresserras!79:79)"
$ws.Range("E59").Formula = "=LOOKUP(2,1/(TransportsTresserras!78:78<>""""),TransportsTresserras!78:78)"
$ws.Range("C60").Formula = "=LOOKUP(2,1/(TransportsTresserras!81:81<>""""),TransportsTresserras!81:81)"
$ws.Range("D60").Formula = "=LOOKUP(2,1/(TransportsTresserras!83:83<>""""),TransportsTresserras!83:83)"
$ws.Range("E60").Formula = "=LOOKUP(2,1/(TransportsTresserras!82:82<>""""),TransportsTresserras!82:82)"
$ws.Range("C61").Formula = "=LOOKUP(2,1/(TransportsTresserras!85:85<>""""),TransportsTresserras!85:85)"
$ws.Range("D61").Formula = "=LOOKUP(2,1/(TransportsTresserras!87:87<>""""),TransportsTresserras!87:87)"
$ws.Range("E61").Formula = "=LOOKUP(2,1/(TransportsTresserras!86:86<>""""),TransportsTresserras!86:86)"
$ws.Range("C62").Formula = "=LOOKUP(2,1/(TransportsTresserras!89:89<>""""),TransportsTresserras!89:89)"
$ws.Range("D62").Formula = "=LOOKUP(2,1/(TransportsTresserras!91:91<>""""),TransportsTresserras!91:91)"
$ws.Range("E62").Formula = "=LOOKUP(2,1/(TransportsTresserras!90:90<>""""),TransportsTresserras!90:90)"
$ws.Range("C63").Formula = "=LOOKUP(2,1/(TransportsTresserras!93:93<>""""),TransportsTresserras!93:93)"
$ws.Range("D63").Formula = "=LOOKUP(2,1/(TransportsTresserras!95:95<>""""),TransportsTresserras!95:95)"
$ws.Range("E63").Formula = "=LOOKUP(2,1/(TransportsTresserras!94:94<>""""),TransportsTresserras!94:94)"
$ws.Range("C64").Formula = "=LOOKUP(2,1/(TransportsTresserras!97:97<>""""),TransportsTresserras!97:97)"
$ws.Range("D64").Formula = "=LOOKUP(2,1/(TransportsTresserras!99:99<>""""),TransportsTresserras!99:99)"
$ws.Range("E64").Formula = "=LOOKUP(2,1/(TransportsTresserras!98:98<>""""),TransportsTresserras!98:98)"
$ws.Range("C65").Formula = "=LOOKUP(2,1/(TransportsTresserras!101:101<>""""),TransportsTresserras!101:101)"
$ws.Range("D65").Formula = "=LOOKUP(2,1/(TransportsTresserras!103:103<>""""),TransportsTresserras!103:103)"
$ws.Range("E65").Formula = "=LOOKUP(2,1/(TransportsTresserras!102:102<>""""),TransportsTresserras!102:102)"
$ws.Range("C66").Formula = "=LOOKUP(2,1/(TransportsTresserras!105:105<>""""),TransportsTresserras!105:105)"
$ws.Range("D66").Formula = "=LOOKUP(2,1/(TransportsTresserras!107:107<>""""),TransportsTresserras!107:107)"
$ws.Range("E66").Formula = "=LOOKUP(2,1/(TransportsTresserras!106:106<>""""),TransportsTresserras!106:106)"

# --- Remove obsolete trailing rows (old TorinoLocal tail, rows 67-70) ---
$ws.Rows("67:70").Delete()

Write-Host "Edit complete"
